$p = $ppt.ActivePresentation

# Slide 10: "Aufbau des GUI (1/2)" -> "Komponenten der User Experience (1/2)"
$s10 = $p.Slides.Item(10)
$s10.Shapes.Item(1).TextFrame.TextRange.Text = "Komponenten der User Experience (1/2)"

# Slide 11: "Aufbau des GUI (2/2)" -> two runs: "Komponenten der User Experience " + "(2/2)"
$s11 = $p.Slides.Item(11)
$tr11 = $s11.Shapes.Item(1).TextFrame.TextRange
$tr11.Text = "Komponenten der User Experience "
$tr11.InsertAfter("(2/2)") | Out-Null

# Slide 4: resize/reposition "Rechteck 9" shape
$s4 = $p.Slides.Item(4)
$rect = $s4.Shapes.Item(6)
$rect.Left = 273.7675690551181
$rect.Width = 408.5581302362205
